$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.3347771808043092
$ws.Range("C2").Value = 0.06664622201573422
$ws.Range("D2").Value = 0.6286083549714903
$ws.Range("E2").Value = 0.2437140734033925
$ws.Range("G2").Value = 0.0025156388566126
$ws.Range("I2").Value = 1.09418501108447
$ws.Range("J2").Value = 0.1173822175759014
$ws.Range("K2").Value = 0.3860263727944755
$ws.Range("M2").Value = 0.3173841013454535
$ws.Range("O2").Value = 5.486096691256932
$ws.Range("B3").Value = 0.3055954685944187
$ws.Range("C3").Value = 0.06140566285115767
$ws.Range("D3").Value = 0.6233357810502014
$ws.Range("E3").Value = 0.2423824831198758
$ws.Range("G3").Value = 0.002518513636595854
$ws.Range("I3").Value = 1.100388261393931
$ws.Range("J3").Value = 0.1172214461334704
$ws.Range("K3").Value = 0.3531534470905342
$ws.Range("M3").Value = 0.3066533222910621
$ws.Range("O3").Value = 5.505245606221024
$ws.Range("B4").Value = 0.2877671563255149
$ws.Range("C4").Value = 0.0582083442036776
$ws.Range("D4").Value = 0.6203888416040968
$ws.Range("E4").Value = 0.2416745891890528
$ws.Range("G4").Value = 0.002520373451876947
$ws.Range("I4").Value = 1.104666085569491
$ws.Range("J4").Value = 0.1171727427703217
$ws.Range("K4").Value = 0.3330760061856068
$ws.Range("M4").Value = 0.3002129441424799
$ws.Range("O4").Value = 5.519472982738193
$ws.Range("B5").Value = 0.2805248495409103
$ws.Range("C5").Value = 0.05691058680129402
$ws.Range("D5").Value = 0.6192610937940231
$ws.Range("E5").Value = 0.2414137395966272
$ws.Range("G5").Value = 0.002521155223481065
$ws.Range("I5").Value = 1.106527280031329
$ws.Range("J5").Value = 0.1171654818358618
$ws.Range("K5").Value = 0.3249215231173679
$ws.Range("M5").Value = 0.29762589207904
$ws.Range("O5").Value = 5.525891935545218
$ws.Range("B6").Value = 0.2793236634665277
$ws.Range("C6").Value = 0.05669540911299009
$ws.Range("D6").Value = 0.6190782539890165
$ws.Range("E6").Value = 0.2413720951420295
$ws.Range("G6").Value = 0.002521286480664425
$ws.Range("I6").Value = 1.106843453979831
$ws.Range("J6").Value = 0.1171650366404577
$ws.Range("K6").Value = 0.3235691347304623
$ws.Range("M6").Value = 0.2971985805532853
$ws.Range("O6").Value = 5.526995319019136
$ws.Range("B7").Value = 0.2876693908889081
$ws.Range("C7").Value = 0.05819082117618279
$ws.Range("D7").Value = 0.620373336044679
$ws.Range("E7").Value = 0.2416709593878288
$ws.Range("G7").Value = 0.002520383898278869
$ws.Range("I7").Value = 1.104690708715317
$ws.Range("J7").Value = 0.1171725938723043
$ws.Range("K7").Value = 0.3329659211415645
$ws.Range("M7").Value = 0.300177902406503
$ws.Range("O7").Value = 5.519557035701581
$ws.Range("B8").Value = 0.3246970065718529
$ws.Range("C8").Value = 0.06483506696672237
$ws.Range("D8").Value = 0.6267301588840581
$ws.Range("E8").Value = 0.2432321917773734
$ws.Range("G8").Value = 0.002516610471490993
$ws.Range("I8").Value = 1.096226556442954
$ws.Range("J8").Value = 0.117316409730325
$ws.Range("K8").Value = 0.3746698863588733
$ws.Range("M8").Value = 0.3136534198836145
$ws.Range("O8").Value = 5.492186677075978
$ws.Range("B9").Value = 0.39800346823165
$ws.Range("C9").Value = 0.07802507391470215
$ws.Range("D9").Value = 0.6414964871345319
$ws.Range("E9").Value = 0.2471630979188717
$ws.Range("G9").Value = 0.002509958822926052
$ws.Range("I9").Value = 1.083349599265802
$ws.Range("J9").Value = 0.1179949327673668
$ws.Range("K9").Value = 0.4572837572846709
$ws.Range("M9").Value = 0.3412512985656946
$ws.Range("O9").Value = 5.458109548560685
$ws.Range("B10").Value = 0.4522728828555103
$ws.Range("C10").Value = 0.08781308872950433
$ws.Range("D10").Value = 0.6537445753925795
$ws.Range("E10").Value = 0.2505802004982556
$ws.Range("G10").Value = 0.002505523277153952
$ws.Range("I10").Value = 1.076157902206148
$ws.Range("J10").Value = 0.1187349525521526
$ws.Range("K10").Value = 0.5184755319343139
$ws.Range("M10").Value = 0.3622382979974432
$ws.Range("O10").Value = 5.445023366289007
$ws.Range("B11").Value = 0.4770483869478426
$ws.Range("C11").Value = 0.09228697978296907
$ws.Range("D11").Value = 0.6596198314249477
$ws.Range("E11").Value = 0.2522494633514967
$ws.Range("G11").Value = 0.002503602475840044
$ws.Range("I11").Value = 1.073379061192711
$ws.Range("J11").Value = 0.119124019732233
$ws.Range("K11").Value = 0.5464186644987876
$ws.Range("M11").Value = 0.3719395156033372
$ws.Range("O11").Value = 5.441666332340816
$ws.Range("B12").Value = 0.4864425698397099
$ws.Range("C12").Value = 0.0939841566910502
$ws.Range("D12").Value = 0.6618881991916226
$ws.Range("E12").Value = 0.2528980536307159
$ws.Range("G12").Value = 0.002502888985733381
$ws.Range("I12").Value = 1.072397644239452
$ws.Range("J12").Value = 0.1192788822088531
$ws.Range("K12").Value = 0.5570150206284552
$ws.Range("M12").Value = 0.3756351675509322
$ws.Range("O12").Value = 5.440768430896071
$ws.Range("B13").Value = 0.4844188270106713
$ws.Range("C13").Value = 0.09361850609781186
$ws.Range("D13").Value = 0.6613977311325527
$ws.Range("E13").Value = 0.2527576358179715
$ws.Range("G13").Value = 0.00250304203255807
$ws.Range("I13").Value = 1.072605857656022
$ws.Range("J13").Value = 0.1192451949790367
$ws.Range("K13").Value = 0.5547322501730889
$ws.Range("M13").Value = 0.3748382662296734
$ws.Range("O13").Value = 5.44094520519991
$ws.Range("B14").Value = 0.477821009220861
$ws.Range("C14").Value = 0.09242654746583412
$ws.Range("D14").Value = 0.6598055795286086
$ws.Range("E14").Value = 0.2523024932314897
$ws.Range("G14").Value = 0.002503543498830632
$ws.Range("I14").Value = 1.073296898974235
$ws.Range("J14").Value = 0.119136609470786
$ws.Range("K14").Value = 0.5472901370426371
$ws.Range("M14").Value = 0.3722431188399398
$ws.Range("O14").Value = 5.441584979024753
$ws.Range("B15").Value = 0.4737812362927514
$ws.Range("C15").Value = 0.09169682915839417
$ws.Range("D15").Value = 0.6588360070042825
$ws.Range("E15").Value = 0.2520258500475592
$ws.Range("G15").Value = 0.002503852466563915
$ws.Range("I15").Value = 1.073729411690685
$ws.Range("J15").Value = 0.1190710782577895
$ws.Range("K15").Value = 0.5427335549401278
$ws.Range("M15").Value = 0.3706563781661103
$ws.Range("O15").Value = 5.442025479402446
$ws.Range("B16").Value = 0.45065548107155
$ws.Range("C16").Value = 0.08752113342842449
$ws.Range("D16").Value = 0.653366711935405
$ws.Range("E16").Value = 0.2504734176702286
$ws.Range("G16").Value = 0.002505650751256638
$ws.Range("I16").Value = 1.076349422087056
$ws.Range("J16").Value = 0.1187105804089796
$ws.Range("K16").Value = 0.5166514913895242
$ws.Range("M16").Value = 0.3616073890399321
$ws.Range("O16").Value = 5.445294996779438
$ws.Range("B17").Value = 0.4364908323653935
$ws.Range("C17").Value = 0.08496489406843466
$ws.Range("D17").Value = 0.650089139944356
$ws.Range("E17").Value = 0.2495504322113504
$ws.Range("G17").Value = 0.002506778725178793
$ws.Range("I17").Value = 1.078082919852264
$ws.Range("J17").Value = 0.1185028491528612
$ws.Range("K17").Value = 0.5006779993445321
$ws.Range("M17").Value = 0.3560954945872936
$ws.Range("O17").Value = 5.447965650730112
$ws.Range("B18").Value = 0.4283520165797086
$ws.Range("C18").Value = 0.08349661803757158
$ws.Range("D18").Value = 0.6482325445047081
$ws.Range("E18").Value = 0.2490303633720998
$ws.Range("G18").Value = 0.002507436635486251
$ws.Range("I18").Value = 1.07912635573971
$ws.Range("J18").Value = 0.1183883029177579
$ws.Range("K18").Value = 0.4915005374843702
$ws.Range("M18").Value = 0.3529397148223836
$ws.Range("O18").Value = 5.449746079381839
$ws.Range("B19").Value = 0.4255977946661176
$ws.Range("C19").Value = 0.08299983086459406
$ws.Range("D19").Value = 0.6476088448041253
$ws.Range("E19").Value = 0.2488561342653952
$ws.Range("G19").Value = 0.002507660962437693
$ws.Range("I19").Value = 1.079487609336468
$ws.Range("J19").Value = 0.1183503674151893
$ws.Range("K19").Value = 0.4883949493662385
$ws.Range("M19").Value = 0.3518737198119339
$ws.Range("O19").Value = 5.450390866175042
$ws.Range("B20").Value = 0.4379978258835422
$ws.Range("C20").Value = 0.08523680292408642
$ws.Range("D20").Value = 0.6504350862551007
$ws.Range("E20").Value = 0.2496475672375595
$ws.Range("G20").Value = 0.002506657706036672
$ws.Range("I20").Value = 1.07789358638955
$ws.Range("J20").Value = 0.1185244517740784
$ws.Range("K20").Value = 0.5023773660303732
$ws.Range("M20").Value = 0.3566807438970869
$ws.Range("O20").Value = 5.44765606656145
$ws.Range("B21").Value = 0.4797586195258248
$ws.Range("C21").Value = 0.09277657325726807
$ws.Range("D21").Value = 0.6602720523968912
$ws.Range("E21").Value = 0.2524357327267737
$ws.Range("G21").Value = 0.002503395829970578
$ws.Range("I21").Value = 1.073091999836635
$ws.Range("J21").Value = 0.1191682993359962
$ws.Range("K21").Value = 0.5494756652558408
$ws.Range("M21").Value = 0.3730047803602119
$ws.Range("O21").Value = 5.441386929411664
$ws.Range("B22").Value = 0.5071228427804613
$ws.Range("C22").Value = 0.09772178026059919
$ws.Range("D22").Value = 0.666954786317433
$ws.Range("E22").Value = 0.2543539837582571
$ws.Range("G22").Value = 0.002501344853788459
$ws.Range("I22").Value = 1.070366978768874
$ws.Range("J22").Value = 0.1196329847597468
$ws.Range("K22").Value = 0.5803437607908393
$ws.Range("M22").Value = 0.3838017109604053
$ws.Range("O22").Value = 5.439465785360142
$ws.Range("B23").Value = 0.4925116742747946
$ws.Range("C23").Value = 0.09508084177537057
$ws.Range("D23").Value = 0.6633649083622402
$ws.Range("E23").Value = 0.2533214015403971
$ws.Range("G23").Value = 0.002502432122132683
$ws.Range("I23").Value = 1.071783568898773
$ws.Range("J23").Value = 0.1193809597332205
$ws.Range("K23").Value = 0.5638610961230768
$ws.Range("M23").Value = 0.378027502949088
$ws.Range("O23").Value = 5.440292008595691
$ws.Range("B24").Value = 0.4373164996625292
$ws.Range("C24").Value = 0.08511386875265714
$ws.Range("D24").Value = 0.6502785975630161
$ws.Range("E24").Value = 0.2496036195726887
$ws.Range("G24").Value = 0.002506712389513582
$ws.Range("I24").Value = 1.077979038189014
$ws.Range("J24").Value = 0.1185146700240622
$ws.Range("K24").Value = 0.5016090639059314
$ws.Range("M24").Value = 0.3564161119252702
$ws.Range("O24").Value = 5.447795266282583
$ws.Range("B25").Value = 0.3780990086581255
$ws.Range("C25").Value = 0.07443969394674355
$ws.Range("D25").Value = 0.6372558443663081
$ws.Range("E25").Value = 0.2460067093051315
$ws.Range("G25").Value = 0.002511678665681184
$ws.Range("I25").Value = 1.086434730374492
$ws.Range("J25").Value = 0.1177689480038531
$ws.Range("K25").Value = 0.4572837572846709
$ws.Range("M25").Value = 0.3336602305542975
$ws.Range("O25").Value = 5.465230015596035

Write-Host "Updated B2:O25 with new pl_mw results (380 kV case)"
